# Update the "想去人数" (number of interested people) column (F) on the
# "展览" and "全部类型" sheets to reflect freshly generated output values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9739
$ws1.Range("F3").Value = 215
$ws1.Range("F5").Value = 559
$ws1.Range("F6").Value = 470

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9739
$ws4.Range("F3").Value = 215
$ws4.Range("F5").Value = 559
$ws4.Range("F7").Value = 470
